$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2, 3, 4, 6, 7
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -3
